# ---------------------------------------------------------------------------
# 0029-oceanic_fracture_zones_a_better_interpretation
# "Update new files add back 0023"
#
# 1) Refresh the auto "Date" placeholder (field type datetimeFigureOut) that
#    lives on the slide master and every slide layout from 7/24/20 -> 7/9/21.
# 2) On slide 2, the braille transcription textbox ("Rectangle 4") gets its
#    wording adjusted:
#       - paragraph 1 absorbs what used to be paragraph 2's text
#       - the old paragraph 3 ("...mid-ocean") is replaced with a shorter line
#       - the remaining two paragraphs stay the same, just shift up
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Date placeholder refresh -------------------------------------------
$newDate = "7/9/21"

$m = $p.SlideMaster
$m.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$customLayouts = $m.CustomLayouts
$customLayouts.Item(2).Shapes.Item(3).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(4).Shapes.Item(4).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(5).Shapes.Item(6).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(8).Shapes.Item(4).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(9).Shapes.Item(4).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(10).Shapes.Item(3).TextFrame.TextRange.Text = $newDate
$customLayouts.Item(11).Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# --- 2) Braille caption textbox on slide 2 ----------------------------------
$s2 = $p.Slides.Item(2)
$brailleShape = $s2.Shapes.Item(14)

$line1 = "⠠⠚⠀⠮⠀⠐⠏⠀⠷⠀⠮ ⠕⠉⠂⠝⠊⠉⠀⠋⠗⠁⠉⠞⠥⠗⠑⠀⠵⠐⠕"
$line2 = "⠆⠞⠀⠮⠀⠍⠊⠙⠤⠕⠉1⠝"
$line3 = "⠗⠊⠙⠛⠑⠎⠀⠊⠎⠀⠁⠉⠞⠊⠧⠑"
$line4 = "⠐⠣⠁⠀⠋⠁⠥⠇⠞⠐⠜"

$brailleShape.TextFrame.TextRange.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4
